$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Delete row 35 ("BannedPaths" rule), shifting rows 36-134 up to 35-133
$ws.Rows.Item(35).Delete()

# Insert a fresh blank row at row 40, shifting rows 40-133 back down to 41-134
# (restoring rows 41+ to their original, unshifted position)
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the "BannedPath" rule
$ws.Cells.Item(40, 1).Value = "BannedPath"
$ws.Cells.Item(40, 2).Value = "Customer packages should not install content under /libs"
$ws.Cells.Item(40, 3).Value = "Bug"
$ws.Cells.Item(40, 4).Value = "Critical"

# Update the active cell selection shown in the sheet view
$ws.Range("A37").Select()
